$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.31267000400635
$ws.Range("C2").Value = 0.3461160285025073
$ws.Range("E2").Value = 0.1097307580660658
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.2883040419709957
$ws.Range("H2").Value = 0.4758576272819894
$ws.Range("L2").Value = 0.2100670150863522
$ws.Range("O2").Value = 1.453752698736565
$ws.Range("B3").Value = 1.167875531889081
$ws.Range("C3").Value = 0.3332947800523414
$ws.Range("E3").Value = 0.1114209922514373
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.2947535778493346
$ws.Range("H3").Value = 0.483456961853598
$ws.Range("L3").Value = 0.1992484599809217
$ws.Range("O3").Value = 1.483065920702316
$ws.Range("B4").Value = 1.078733597952464
$ws.Range("C4").Value = 0.325481867613945
$ws.Range("E4").Value = 0.1125408572756683
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.2991058866762124
$ws.Range("H4").Value = 0.488454943375686
$ws.Range("L4").Value = 0.1926865281558037
$ws.Range("O4").Value = 1.502581903982275
$ws.Range("B5").Value = 1.042350318477077
$ws.Range("C5").Value = 0.3223133461048917
$ws.Range("E5").Value = 0.1130178338913783
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.3009777669265716
$ws.Range("H5").Value = 0.4905750567884546
$ws.Range("L5").Value = 0.1900329119019943
$ws.Range("O5").Value = 1.510915492073046
$ws.Range("B6").Value = 1.036305521290103
$ws.Range("C6").Value = 0.3217881508127221
$ws.Range("E6").Value = 0.1130982807282876
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.3012945164516658
$ws.Range("H6").Value = 0.4909321345728017
$ws.Range("L6").Value = 0.1895935183806046
$ws.Range("O6").Value = 1.512322241920039
$ws.Range("B7").Value = 1.078243148643935
$ws.Range("C7").Value = 0.3254390733669936
$ws.Range("E7").Value = 0.1125472064634394
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.2991307340576626
$ws.Range("H7").Value = 0.4884831984069962
$ws.Range("L7").Value = 0.1926506576877642
$ws.Range("O7").Value = 1.502692753602432
$ws.Range("B8").Value = 1.262795629570974
$ws.Range("C8").Value = 0.3416831531455671
$ws.Range("E8").Value = 0.1102965141189713
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.2904462113639852
$ws.Range("H8").Value = 0.4784089332180059
$ws.Range("L8").Value = 0.2063201012503981
$ws.Range("O8").Value = 1.463544373522694
$ws.Range("B9").Value = 1.622725288564595
$ws.Range("C9").Value = 0.3739945409050165
$ws.Range("E9").Value = 0.1065344363646581
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.2765443622809158
$ws.Range("H9").Value = 0.4612906122041593
$ws.Range("L9").Value = 0.2337620226353181
$ws.Range("O9").Value = 1.398856741716443
$ws.Range("B10").Value = 1.885862654205766
$ws.Range("C10").Value = 0.3979955941499043
$ws.Range("E10").Value = 0.1041682225330671
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.2682596985447532
$ws.Range("H10").Value = 0.4503258568457653
$ws.Range("L10").Value = 0.2543081269858192
$ws.Range("O10").Value = 1.358752058093231
$ws.Range("B11").Value = 2.005269001371573
$ws.Range("C11").Value = 0.408967877153799
$ws.Range("E11").Value = 0.103178272590922
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.264914306388718
$ws.Range("H11").Value = 0.4456886690764605
$ws.Range("L11").Value = 0.263738029756297
$ws.Range("O11").Value = 1.342130716118191
$ws.Range("B12").Value = 2.050440376720587
$ws.Range("C12").Value = 0.4131302494783711
$ws.Range("E12").Value = 0.1028158480586114
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.2637087512102667
$ws.Range("H12").Value = 0.4439832073610717
$ws.Range("L12").Value = 0.2673207821267454
$ws.Range("O12").Value = 1.336070950738801
$ws.Range("B13").Value = 2.040713965452596
$ws.Range("C13").Value = 0.4122334845456805
$ws.Range("E13").Value = 0.1028933487802046
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.263965657705711
$ws.Range("H13").Value = 0.4443482595714272
$ws.Range("L13").Value = 0.2665486473097474
$ws.Range("O13").Value = 1.337365588786511
$ws.Range("B14").Value = 2.008986198792854
$ws.Range("C14").Value = 0.4093101716894125
$ws.Range("E14").Value = 0.1031482061369893
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.2648138943912457
$ws.Range("H14").Value = 0.445547346201856
$ws.Range("L14").Value = 0.2640325481815182
$ws.Range("O14").Value = 1.341627472242948
$ws.Range("B15").Value = 1.989546059411907
$ws.Range("C15").Value = 0.4075205118413123
$ws.Range("E15").Value = 0.1033059352068051
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.2653414550234103
$ws.Range("H15").Value = 0.4462884063953254
$ws.Range("L15").Value = 0.2624929043619346
$ws.Range("O15").Value = 1.344268553138534
$ws.Range("B16").Value = 1.878052983760426
$ws.Range("C16").Value = 0.3972795861907059
$ws.Range("E16").Value = 0.104234659130011
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.2684868804338478
$ws.Range("H16").Value = 0.4506359746689697
$ws.Range("L16").Value = 0.2536935267214631
$ws.Range("O16").Value = 1.359871041522268
$ws.Range("B17").Value = 1.809577860557852
$ws.Range("C17").Value = 0.3910106971393077
$ws.Range("E17").Value = 0.1048265553746894
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.2705252077755702
$ws.Range("H17").Value = 0.4533929736443554
$ws.Range("L17").Value = 0.2483166509104819
$ws.Range("O17").Value = 1.369858978547271
$ws.Range("B18").Value = 1.770165060356362
$ws.Range("C18").Value = 0.387410110389709
$ws.Range("E18").Value = 0.1051751348073591
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.2717374232494478
$ws.Range("H18").Value = 0.4550117330351995
$ws.Range("L18").Value = 0.2452318701056555
$ws.Range("O18").Value = 1.37575642785886
$ws.Range("B19").Value = 1.756815882931335
$ws.Range("C19").Value = 0.3861919030056526
$ws.Range("E19").Value = 0.1052945547449742
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.2721546862038764
$ws.Range("H19").Value = 0.4555654828311617
$ws.Range("L19").Value = 0.2441887697705596
$ws.Range("O19").Value = 1.377779388463182
$ws.Range("B20").Value = 1.816870042405355
$ws.Range("C20").Value = 0.391677504662141
$ws.Range("E20").Value = 0.1047627048073352
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.2703040998918809
$ws.Range("H20").Value = 0.4530960695821307
$ws.Range("L20").Value = 0.2488882164510358
$ws.Range("O20").Value = 1.368779938943234
$ws.Range("B21").Value = 2.018306660957478
$ws.Range("C21").Value = 0.4101686214441429
$ws.Range("E21").Value = 0.103073010389032
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.2645630806586325
$ws.Range("H21").Value = 0.4451937728944486
$ws.Range("L21").Value = 0.2647712668891558
$ws.Range("O21").Value = 1.340369284396928
$ws.Range("B22").Value = 2.149692374691369
$ws.Range("C22").Value = 0.4222966396738457
$ws.Range("E22").Value = 0.1020412614511397
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.2611683164326735
$ws.Range("H22").Value = 0.4403238034972361
$ws.Range("L22").Value = 0.2752207830668993
$ws.Range("O22").Value = 1.323167893303079
$ws.Range("B23").Value = 2.079594434594526
$ws.Range("C23").Value = 0.4158198743010928
$ws.Range("E23").Value = 0.1025852803375198
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.2629473405389788
$ws.Range("H23").Value = 0.4428960036351199
$ws.Range("L23").Value = 0.2696374088361182
$ws.Range("O23").Value = 1.33222320349438
$ws.Range("B24").Value = 1.813573388645409
$ws.Range("C24").Value = 0.3913760301095692
$ws.Range("E24").Value = 0.1047915458263784
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.2704039370964892
$ws.Range("H24").Value = 0.4532301948680555
$ws.Range("L24").Value = 0.2486297915496323
$ws.Range("O24").Value = 1.369267289230834
$ws.Range("B25").Value = 1.525576640175245
$ws.Range("C25").Value = 0.365206285594752
$ws.Range("E25").Value = 0.1074823655706556
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.2799679931007475
$ws.Range("H25").Value = 0.4656387851305581
$ws.Range("L25").Value = 0.2262704932107624
$ws.Range("O25").Value = 1.415057202656271
